$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A96's date/time value (R script result re-run gave an updated timestamp)
$ws.Range("A96").Value = 45475.2916666667

# Copy the date formatting from A96 onto the new A97 cell before filling values,
# so the new row's date cell keeps the same number format / style as the rest
# of column A.
$ws.Range("A96").Copy()
$ws.Range("A97").PasteSpecial(-4122)

# Append the new data row (row 97) pulled in by the latest R script run.
$ws.Range("A97").Value = 45476.6359722222
$ws.Range("B97").Value = 20400
$ws.Range("C97").Value = 6.23999977111816
$ws.Range("D97").Value = 6
$ws.Range("E97").Value = 6.15999984741211
$ws.Range("F97").Value = 6

# adj_close (G97) came back from R as the character "6" rather than a plain
# number, so write it as text. A bare Value = "6" would be auto-coerced to a
# number by Excel's type sniffer, so build it as a formula that evaluates to
# the text string "6" on a scratch cell, then paste-special just the
# resulting value (not the formula, not the formatting) onto G97 - that
# locks it in as literal text - and tidy up the scratch cell afterward.
$ws.Range("ZZ1000").Formula = "=""6"""
$ws.Range("ZZ1000").Copy()
$ws.Range("G97").PasteSpecial(-4163)
$ws.Range("ZZ1000").Clear()

$ws.Range("H97").Value = "PAL.MI"
